$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.212.70'
$ws.Range("E2").Value = '  +3.22%  '
$ws.Range("D3").Value = '2.543.42'
$ws.Range("E3").Value = '  +3.45%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''541.13'
$ws.Range("E5").Value = '  +2.36%  '
$ws.Range("D6").Value = '''145.44'
$ws.Range("E6").Value = '  +2.80%  '
$ws.Range("D7").Value = '''0.995'
$ws.Range("E7").Value = '  -0.30%  '
$ws.Range("E8").Value = '  +1.58%  '
$ws.Range("D9").Value = '2.572.37'
$ws.Range("E9").Value = '  +4.03%  '
$ws.Range("E10").Value = '  +3.74%  '
$ws.Range("E11").Value = '  +1.81%  '
$ws.Range("D12").Value = '''5.47'
$ws.Range("E12").Value = '  -0.35%  '
$ws.Range("E13").Value = '  +5.46%  '
$ws.Range("D14").Value = '2.994.29'
$ws.Range("E14").Value = '  +3.63%  '
$ws.Range("D15").Value = '''24.29'
$ws.Range("E15").Value = '  +4.31%  '
$ws.Range("D16").Value = '60.154.55'
$ws.Range("E16").Value = '  +3.34%  '
$ws.Range("E17").Value = '  +6.40%  '
$ws.Range("D18").Value = '2.569.80'
$ws.Range("E18").Value = '  +4.15%  '
$ws.Range("D19").Value = '''11.30'
$ws.Range("E19").Value = '  +1.76%  '
$ws.Range("E20").Value = '  +3.59%  '
$ws.Range("D21").Value = '''327.98'
$ws.Range("E21").Value = '  +3.21%  '
$ws.Range("D22").Value = '''0.998'
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("E23").Value = '  +5.34%  '
$ws.Range("D24").Value = '''63.10'
$ws.Range("E24").Value = '  +5.26%  '
$ws.Range("E25").Value = '  +2.40%  '
$ws.Range("D26").Value = '''0.167'
$ws.Range("E26").Value = '  +5.51%  '
$ws.Range("D27").Value = '''0.994'
$ws.Range("E27").Value = '  -0.46%  '
$ws.Range("D28").Value = '''8.02'
$ws.Range("E28").Value = '  +5.71%  '
$ws.Range("D29").Value = '''7.15'
$ws.Range("E29").Value = '  +7.53%  '
$ws.Range("D30").Value = '0.0₃0802'
$ws.Range("E30").Value = '  +7.73%  '
$ws.Range("D31").Value = '''1.81'
$ws.Range("E31").Value = '  +3.44%  '
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("D33").Value = '''165.18'
$ws.Range("E33").Value = '  +5.45%  '
$ws.Range("D34").Value = '''1.49'
$ws.Range("E34").Value = '  +7.59%  '
$ws.Range("D35").Value = '''0.998'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = '''18.80'
$ws.Range("E36").Value = '  +2.79%  '
$ws.Range("D37").Value = '''4.46'
$ws.Range("E37").Value = '  +5.07%  '
$ws.Range("E38").Value = '  +5.19%  '
$ws.Range("D39").Value = '''37.05'
$ws.Range("E39").Value = '  +1.51%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = '''303.66'
$ws.Range("E40").Value = '  +2.53%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '''5.61'
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("D42").Value = '''0.843'
$ws.Range("E42").Value = '  +9.83%  '
$ws.Range("E43").Value = '  +4.81%  '
$ws.Range("D44").Value = '''0.610'
$ws.Range("E44").Value = '  +4.34%  '
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("E46").Value = '  +1.05%  '
$ws.Range("D47").Value = '''127.13'
$ws.Range("E47").Value = '  +3.44%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '''0.0941'
$ws.Range("E48").Value = '  +2.93%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''19.02'
$ws.Range("E49").Value = '  +4.57%  '
$ws.Range("E50").Value = '  +3.42%  '
$ws.Range("D51").Value = '''0.0230'
$ws.Range("E51").Value = '  +3.19%  '
